$d = $word.ActiveDocument

# Update date paragraph
$d.Paragraphs.Item(1).Range.Text = "2025-09-07 Sunday"

# Update table cells (row-major order, 20 rows x 5 cols)
$t = $d.Tables.Item(1)
$values = @(
  "61-15=",
  "11+47=",
  "21+59=",
  "85-49=",
  "46-5=",
  "49+50=",
  "95-15=",
  "13+44=",
  "89-56=",
  "6+22=",
  "51+17=",
  "15+37=",
  "64+19=",
  "96-18=",
  "76-33=",
  "74-54=",
  "60-33=",
  "82-5=",
  "38+3=",
  "22+66=",
  "62+26=",
  "76-65=",
  "26+58=",
  "18+13=",
  "34+38=",
  "73+0=",
  "1+85=",
  "79-20=",
  "42+2=",
  "33-6=",
  "54+37=",
  "94-76=",
  "69-39=",
  "62-45=",
  "85-63=",
  "80-36=",
  "40+22=",
  "70-65=",
  "6+48=",
  "3+39=",
  "48+20=",
  "3+2=",
  "70-42=",
  "84-82=",
  "51-21=",
  "76-75=",
  "5+1=",
  "45+14=",
  "13+65=",
  "13+80=",
  "66+5=",
  "17+23=",
  "7+5=",
  "36+41=",
  "44+10=",
  "49+22=",
  "40-31=",
  "11-6=",
  "89-40=",
  "64+6=",
  "19-17=",
  "78-37=",
  "16+19=",
  "47+33=",
  "58-53=",
  "96+1=",
  "30+29=",
  "59-8=",
  "68+11=",
  "38+14=",
  "37+6=",
  "91-43=",
  "54+21=",
  "54-16=",
  "26+7=",
  "55-30=",
  "35+47=",
  "82-47=",
  "70-8=",
  "44+52=",
  "90-64=",
  "71-58=",
  "57-28=",
  "15-4=",
  "38+5=",
  "62+23=",
  "54-39=",
  "38+23=",
  "84-69=",
  "76-28=",
  "10+83=",
  "33-11=",
  "43-5=",
  "54+34=",
  "34+59=",
  "31+68=",
  "44+53=",
  "65+30=",
  "13+62=",
  "64-2="
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
  for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $values[$idx]
    $idx = $idx + 1
  }
}

Write-Host "Replaced" $idx "cells"